# Fruta / hortaliza, semanal
#
# Insert a new weekly data row at row 484 (pushing the existing rows 484-553
# down to 485-554) and populate it with the new observation's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 484; existing rows 484..553 shift to 485..554.
$ws.Rows.Item(484).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(484, 1).Value2  = 8
$ws.Cells.Item(484, 2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(484, 3).Value2  = "Coquimbo"
$ws.Cells.Item(484, 4).Value2  = 45131
$ws.Cells.Item(484, 5).Value2  = 4
$ws.Cells.Item(484, 6).Value2  = 100112032
$ws.Cells.Item(484, 7).Value2  = "Zapallo italiano"
$ws.Cells.Item(484, 8).Value2  = "Sin especificar"
$ws.Cells.Item(484, 9).Value2  = "Primera"
$ws.Cells.Item(484, 10).Value2 = 500
$ws.Cells.Item(484, 11).Value2 = 15000
$ws.Cells.Item(484, 12).Value2 = 16000
$ws.Cells.Item(484, 13).Value2 = 15500
$ws.Cells.Item(484, 14).Value2 = "$/caja 50 unidades"
$ws.Cells.Item(484, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(484, 16).Value2 = 310
$ws.Cells.Item(484, 17).Value2 = 50
$ws.Cells.Item(484, 18).Value2 = "Hortaliza"
